$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first")
$ws.Activate()
try { $excel.ActiveWindow.TopLeftCell = $ws.Range("A12") ; Write-Output "TopLeftCell prop worked" } catch { Write-Output "TopLeftCell prop failed: $_" }
try { $ws.Range("A1").Application.Goto($ws.Range("A12"), $false) ; Write-Output "Goto worked" } catch { Write-Output "Goto failed: $_" }
